$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the desired effect (power) value in A3; dependent formulas in
# B3 (strøm), E3 (R3) and H3 recalc automatically.
$ws.Range("A3").Value = 5
